$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.723.70'
$ws.Range("E2").Value = '  +0.66%  '
$ws.Range("D3").Value = '3.044.75'
$ws.Range("E3").Value = '  -0.10%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '538.36'
$ws.Range("E5").Value = '  +0.97%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '133.65'
$ws.Range("E6").Value = '  +2.81%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").Value = '3.036.86'
$ws.Range("E8").Value = '  -0.18%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.491'
$ws.Range("E9").Value = '  +1.77%  '
$ws.Range("D10").Value = '0.152'
$ws.Range("E10").Value = '  +1.12%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.15'
$ws.Range("E11").Value = '  +1.00%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.448'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000221'
$ws.Range("E13").Value = '  -0.26%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '33.84'
$ws.Range("E14").Value = '  +0.56%  '
$ws.Range("D15").Value = '3.547.08'
$ws.Range("E15").Value = '  +1.32%  '
$ws.Range("D16").Value = '62.759.68'
$ws.Range("E16").Value = '  +0.74%  '
$ws.Range("E17").Value = '  +1.91%  '
$ws.Range("D18").Value = '3.049.35'
$ws.Range("E18").Value = '  -0.20%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.58'
$ws.Range("E19").Value = '  +1.53%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '466.46'
$ws.Range("E20").Value = '  -1.20%  '
$ws.Range("D21").Value = '13.31'
$ws.Range("E21").Value = '  +1.90%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.685'
$ws.Range("E22").Value = '  -0.29%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.94'
$ws.Range("E23").Value = '  -1.36%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '78.17'
$ws.Range("E24").Value = '  +0.71%  '
$ws.Range("D25").Value = '12.01'
$ws.Range("E25").Value = '  +1.79%  '
$ws.Range("D26").Value = '0.998'
$ws.Range("E26").Value = '  -0.03%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.69'
$ws.Range("E27").Value = '  +1.44%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.73'
$ws.Range("E28").Value = '  -3.45%  '
$ws.Range("E29").Value = '  +0.03%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '25.79'
$ws.Range("E30").Value = '  +1.75%  '
$ws.Range("E31").Value = '  +5.93%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.86'
$ws.Range("E32").Value = '  +0.50%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '57.50'
$ws.Range("E33").Value = '  -0.92%  '
$ws.Range("D34").Value = '2.27'
$ws.Range("E34").Value = '  -3.96%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.39'
$ws.Range("E35").Value = '  +5.33%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.89'
$ws.Range("E36").Value = '  +0.97%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '459.78'
$ws.Range("E37").Value = '  -1.21%  '
$ws.Range("D38").Value = '3.194.06'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0389'
$ws.Range("E39").Value = '  +1.65%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0788'
$ws.Range("E40").Value = '  +1.65%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.116'
$ws.Range("E41").Value = '  +3.53%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.06'
$ws.Range("E42").Value = '  +1.76%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.49'
$ws.Range("E43").Value = '  +0.19%  '
$ws.Range("D45").Value = '0.247'
$ws.Range("E45").Value = '  +1.31%  '
$ws.Range("D46").Value = '25.07'
$ws.Range("E46").Value = '  +4.81%  '
$ws.Range("D47").Value = '121.58'
$ws.Range("E47").Value = '  +4.48%  '
$ws.Range("D48").Value = '0.108'
$ws.Range("E48").Value = '  +2.67%  '
$ws.Range("D49").Value = '1.95'
$ws.Range("E49").Value = '  -1.50%  '
$ws.Range("D50").Value = '0.0₃0511'
$ws.Range("E50").Value = '  +1.78%  '
$ws.Range("E51").Value = '  +7.39%  '
